# Reorder the "Recorded By" (column G) list of names/emails so that
# "System" (or, if absent, "admin@admin.com") is moved to the front of
# the comma-separated list, keeping the relative order of the remaining
# entries unchanged. Single-entry cells and cells already starting with
# the priority entry are left untouched.
#
# NOTE: this runtime's -ceq/-cne/-clike/-cmatch operators behave
# case-INsensitively (confirmed empirically), so exact, case-sensitive
# string matching (needed to tell "System" apart from "system") uses
# [string]::CompareTo(), which is ordinal/case-sensitive and cheap
# (a per-character char-code loop blows the interpreter's statement
# budget across ~150 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $value = $cell.Value2

    if ($null -eq $value) { continue }
    if (-not ($value -is [string])) { continue }
    if ($value -notmatch ',') { continue }

    $parts = @($value -split ',\s*' | ForEach-Object { $_.Trim() })

    # Find the (case-sensitive) index of the priority token: "System"
    # wins over "admin@admin.com" when both would be present.
    $priorityIndex = -1
    for ($i = 0; $i -lt $parts.Count; $i++) {
        if ($parts[$i].CompareTo('System') -eq 0) { $priorityIndex = $i; break }
    }
    if ($priorityIndex -lt 0) {
        for ($i = 0; $i -lt $parts.Count; $i++) {
            if ($parts[$i].CompareTo('admin@admin.com') -eq 0) { $priorityIndex = $i; break }
        }
    }

    if ($priorityIndex -lt 0) { continue }

    $priority = $parts[$priorityIndex]
    $rest = @()
    for ($i = 0; $i -lt $parts.Count; $i++) {
        if ($i -ne $priorityIndex) { $rest += $parts[$i] }
    }

    $newParts = @($priority) + $rest
    $newValue = [string]::Join(', ', $newParts)

    if ($newValue.CompareTo($value) -ne 0) {
        $cell.Value2 = $newValue
    }
}
